$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "day" question rows (33-46) listed one row per weekday with a
# +/- operation baked into the question text. The fix moves that
# variation into the "operands" column (a1:5*) so the question text is
# just "{a}" and the actual day name / operator is generated from the
# sheet at runtime - collapsing 14 near-duplicate rows down to 3.
$ws.Range("A33").Value = "{a}"
$ws.Range("A34").Value = "{a}"
$ws.Range("C35").Value = "a1:5*"

# Remove the now-redundant extra "day" rows (36-46); rows 47+ (currency,
# time, story, distance, addition, subtraction, multiplication, division,
# reminder, percentage) shift up to fill rows 36-45.
$ws.Rows("36:46").Delete()

# Restore the view to where it was left (scrolled down, C35 selected).
$ws.Range("C35").Select()
$excel.ActiveWindow.ScrollRow = 25
